$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vaccine approval")

# Delete rows 15 and 16 (duplicate/leftover rows); remaining rows shift up automatically,
# which fixes the comments (G17->G15, G18->G16, G19->G17), hyperlinks (C17->C15, C18->C16, C19->C17)
# and shared formula ranges.
$ws.Rows.Item(15).Resize(2).Delete() | Out-Null

# Correct the shipping date in (old) row 13, which is unaffected by the deletion since it's above row 15.
$ws.Range("I13").Value = 44230

# Fix the typo in the header text for column M (shared string "shipping_volume_forecast_cummulated" -> "...cumulated")
$ws.Range("M1").Value = "shipping_volume_forecast_cumulated"
